# "new report and excel"
# - Capitalize the PDBI / Additive column headers and make them bold black
# - Add a blank centered formatting block in column K:N (rows 2, 3-13, 14)
# - Hide gridlines and move the active selection to the results table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update header labels (pdbi -> PDBI, additive -> Additive) for all three tables
foreach ($addr in @("C5", "G5", "L5")) {
    $ws.Range($addr).Value = "PDBI"
}
foreach ($addr in @("D5", "H5", "M5")) {
    $ws.Range($addr).Value = "Additive"
}

# 2) Bold + explicit black font color on those same header cells
foreach ($addr in @("C5", "D5", "G5", "H5", "L5", "M5")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Color = 0
}

# 3) New blank formatting block: column N alongside the existing tables, plus
#    the blank border rows above/below the third table (row 2 and row 14)
foreach ($addr in @("K2:N2", "N3:N13", "K14:N14")) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# 4) Turn off gridlines for the sheet view
$excel.ActiveWindow.DisplayGridlines = $false

# 5) Move the selection onto the third (PDBI) table
$ws.Range("K3:M13").Select()

Write-Output "done"
